$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 210.6875
$ws.Cells.Item(2, 9).Value = 223.33333
$ws.Cells.Item(2, 11).Value = 223.33333
$ws.Cells.Item(2, 13).Value = -110.33333
$ws.Cells.Item(12, 8).Value = 200
$ws.Cells.Item(12, 9).Value = 200
$ws.Cells.Item(12, 11).Value = 200
$ws.Cells.Item(12, 13).Value = -30
$ws.Cells.Item(15, 8).Value = 161.48
$ws.Cells.Item(15, 9).Value = 161.48
$ws.Cells.Item(15, 11).Value = 484.4399999999999
$ws.Cells.Item(15, 13).Value = -315.4399999999999
$ws.Cells.Item(17, 8).Value = 1617492.8
$ws.Cells.Item(17, 10).Value = 1617492.8
$ws.Cells.Item(17, 12).Value = 4852478.4
$ws.Cells.Item(17, 14).Value = -4852814.4
$ws.Cells.Item(87, 8).Value = 15400.852
$ws.Cells.Item(87, 10).Value = 15400.852
$ws.Cells.Item(87, 12).Value = 15400.852
$ws.Cells.Item(87, 14).Value = -17896.852
$ws.Cells.Item(88, 8).Value = 3726
$ws.Cells.Item(88, 10).Value = 6502
$ws.Cells.Item(88, 12).Value = 6502
$ws.Cells.Item(88, 14).Value = -7314
$ws.Cells.Item(90, 8).Value = 15400.852
$ws.Cells.Item(90, 10).Value = 15400.852
$ws.Cells.Item(90, 12).Value = 46202.556
$ws.Cells.Item(90, 14).Value = -58682.556
$ws.Cells.Item(91, 8).Value = 3726
$ws.Cells.Item(91, 10).Value = 6502
$ws.Cells.Item(91, 12).Value = 6502
$ws.Cells.Item(91, 14).Value = -9310
$ws.Cells.Item(100, 8).Value = 5778.3335
$ws.Cells.Item(100, 9).Value = 5668.3335
$ws.Cells.Item(100, 10).Value = 5833.3335
$ws.Cells.Item(100, 11).Value = 5668.3335
$ws.Cells.Item(100, 12).Value = 5833.3335
$ws.Cells.Item(100, 13).Value = -5127.3335
$ws.Cells.Item(100, 14).Value = -6915.3335
$ws.Cells.Item(108, 8).Value = 18742.334
$ws.Cells.Item(108, 10).Value = 18742.334
$ws.Cells.Item(108, 12).Value = 18742.334
$ws.Cells.Item(108, 14).Value = -26422.334
$ws.Cells.Item(112, 8).Value = 1241.4445
$ws.Cells.Item(112, 10).Value = 1310.3914
$ws.Cells.Item(112, 12).Value = 3931.1742
$ws.Cells.Item(112, 14).Value = -6147.174199999999
$ws.Cells.Item(132, 8).Value = 13201.549
$ws.Cells.Item(132, 9).Value = 14254.28
$ws.Cells.Item(132, 10).Value = 1922.2858
$ws.Cells.Item(132, 11).Value = 42762.84
$ws.Cells.Item(132, 12).Value = 5766.857400000001
$ws.Cells.Item(132, 13).Value = -40232.84
$ws.Cells.Item(132, 14).Value = -10826.8574
$ws.Cells.Item(137, 8).Value = 1477.9565
$ws.Cells.Item(137, 9).Value = 1777
$ws.Cells.Item(137, 10).Value = 1285.7142
$ws.Cells.Item(137, 11).Value = 5331
$ws.Cells.Item(137, 12).Value = 3857.1426
$ws.Cells.Item(137, 13).Value = -2781
$ws.Cells.Item(137, 14).Value = -8957.142599999999
$ws.Cells.Item(138, 8).Value = 4985.7114
$ws.Cells.Item(138, 9).Value = 3225.8215
$ws.Cells.Item(138, 10).Value = 5699.8696
$ws.Cells.Item(138, 11).Value = 9677.4645
$ws.Cells.Item(138, 12).Value = 17099.6088
$ws.Cells.Item(138, 13).Value = -4537.4645
$ws.Cells.Item(138, 14).Value = -27379.6088

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 400
$ws.Cells.Item(4, 9).Value = 400
$ws.Cells.Item(4, 11).Value = 400
$ws.Cells.Item(4, 13).Value = -284
$ws.Cells.Item(32, 8).Value = 15209
$ws.Cells.Item(32, 9).Value = 12306.381
$ws.Cells.Item(32, 10).Value = 27400
$ws.Cells.Item(32, 11).Value = 12306.381
$ws.Cells.Item(32, 12).Value = 27400
$ws.Cells.Item(32, 13).Value = -12019.381
$ws.Cells.Item(32, 14).Value = -27974
$ws.Cells.Item(41, 8).Value = 841.2857
$ws.Cells.Item(41, 9).Value = 841.2857
$ws.Cells.Item(41, 11).Value = 841.2857
$ws.Cells.Item(41, 13).Value = -427.2857

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 261
$ws.Cells.Item(22, 9).Value = 177.63637
$ws.Cells.Item(22, 11).Value = 177.63637
$ws.Cells.Item(22, 13).Value = -4.636369999999999
$ws.Cells.Item(99, 8).Value = 2382.3845
$ws.Cells.Item(99, 9).Value = 1583.3334
$ws.Cells.Item(99, 11).Value = 1583.3334
$ws.Cells.Item(99, 13).Value = -85.33339999999998
$ws.Cells.Item(107, 8).Value = 524.5
$ws.Cells.Item(107, 9).Value = 566.6667
$ws.Cells.Item(107, 10).Value = 398
$ws.Cells.Item(107, 11).Value = 566.6667
$ws.Cells.Item(107, 12).Value = 398
$ws.Cells.Item(107, 13).Value = 1353.3333
$ws.Cells.Item(107, 14).Value = -4238
$ws.Cells.Item(140, 8).Value = 59763.332
$ws.Cells.Item(140, 10).Value = 59763.332
$ws.Cells.Item(140, 12).Value = 59763.332
$ws.Cells.Item(140, 14).Value = -70123.33199999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2472048.2
$ws.Cells.Item(99, 9).Value = 2911511.5
$ws.Cells.Item(99, 10).Value = 55000
$ws.Cells.Item(99, 11).Value = 2911511.5
$ws.Cells.Item(99, 12).Value = 55000
$ws.Cells.Item(99, 13).Value = -2910013.5
$ws.Cells.Item(99, 14).Value = -57996
$ws.Cells.Item(126, 8).Value = 2472048.2
$ws.Cells.Item(126, 9).Value = 2911511.5
$ws.Cells.Item(126, 10).Value = 55000
$ws.Cells.Item(126, 11).Value = 8734534.5
$ws.Cells.Item(126, 12).Value = 165000
$ws.Cells.Item(126, 13).Value = -8732064.5
$ws.Cells.Item(126, 14).Value = -169940
$ws.Cells.Item(134, 8).Value = 2014.2051
$ws.Cells.Item(134, 9).Value = 2018.6857
$ws.Cells.Item(134, 10).Value = 1975
$ws.Cells.Item(134, 11).Value = 6056.0571
$ws.Cells.Item(134, 12).Value = 5925
$ws.Cells.Item(134, 13).Value = -3521.0571
$ws.Cells.Item(134, 14).Value = -10995
$ws.Cells.Item(138, 8).Value = 40375.445
$ws.Cells.Item(138, 10).Value = 40375.445
$ws.Cells.Item(138, 12).Value = 40375.445
$ws.Cells.Item(138, 14).Value = -50655.445

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1078.7142
$ws.Cells.Item(113, 9).Value = 1179.8125
$ws.Cells.Item(113, 10).Value = 755.2
$ws.Cells.Item(113, 11).Value = 3539.4375
$ws.Cells.Item(113, 12).Value = 2265.6
$ws.Cells.Item(113, 13).Value = -1369.4375
$ws.Cells.Item(113, 14).Value = -6605.6
$ws.Cells.Item(122, 8).Value = 993.2258
$ws.Cells.Item(122, 10).Value = 3724.75
$ws.Cells.Item(122, 12).Value = 33522.75
$ws.Cells.Item(122, 14).Value = -38422.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 1154.1666
$ws.Cells.Item(2, 9).Value = 1046.5
$ws.Cells.Item(2, 11).Value = 1046.5
$ws.Cells.Item(2, 13).Value = -933.5
$ws.Cells.Item(132, 8).Value = 2777.2632
$ws.Cells.Item(132, 9).Value = 2572.6191
$ws.Cells.Item(132, 10).Value = 3030.0588
$ws.Cells.Item(132, 11).Value = 7717.8573
$ws.Cells.Item(132, 12).Value = 9090.1764
$ws.Cells.Item(132, 13).Value = -5187.8573
$ws.Cells.Item(132, 14).Value = -14150.1764

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127, 8).Value = 51163.89
$ws.Cells.Item(127, 10).Value = 51163.89
$ws.Cells.Item(127, 12).Value = 51163.89
$ws.Cells.Item(127, 14).Value = -61083.89
$ws.Cells.Item(133, 8).Value = 23266
$ws.Cells.Item(133, 10).Value = 23266
$ws.Cells.Item(133, 12).Value = 23266
$ws.Cells.Item(133, 14).Value = -28326

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 19950
$ws.Cells.Item(15, 10).Value = 19950
$ws.Cells.Item(15, 12).Value = 19950
$ws.Cells.Item(15, 14).Value = -20526
$ws.Cells.Item(51, 8).Value = 70000
$ws.Cells.Item(51, 10).Value = 70000
$ws.Cells.Item(51, 12).Value = 70000
$ws.Cells.Item(51, 14).Value = -71020
$ws.Cells.Item(107, 8).Value = 765
$ws.Cells.Item(107, 9).Value = 768.3333
$ws.Cells.Item(107, 10).Value = 750
$ws.Cells.Item(107, 11).Value = 2304.9999
$ws.Cells.Item(107, 12).Value = 2250
$ws.Cells.Item(107, 13).Value = -384.9998999999998
$ws.Cells.Item(107, 14).Value = -6090
